$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.812.27"
$ws.Range("E2").Value = "  +6.27%  "
$ws.Range("D3").Value = "3.007.47"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.65"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.83"
$ws.Range("E6").Value = "  +12.39%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").Value = "3.004.26"
$ws.Range("E9").Value = "  +3.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  -4.88%  "
$ws.Range("E11").Value = "  +7.08%  "
$ws.Range("E12").Value = "  +7.39%  "
$ws.Range("E13").Value = "  +8.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.65"
$ws.Range("E14").Value = "  +7.61%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "65.844.15"
$ws.Range("E16").Value = "  +6.40%  "
$ws.Range("D17").Value = "3.508.39"
$ws.Range("E17").Value = "  +3.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.96"
$ws.Range("E18").Value = "  +7.48%  "
$ws.Range("D19").Value = "3.012.08"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "457.63"
$ws.Range("E20").Value = "  +6.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").Value = "  +8.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +5.80%  "
$ws.Range("E23").Value = "  +7.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.34"
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +12.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.39"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  +5.45%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +16.90%  "
$ws.Range("E30").Value = "  +16.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000105"
$ws.Range("E31").Value = "  -5.94%  "
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.96"
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("E37").Value = "  +8.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.16"
$ws.Range("E38").Value = "  +13.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.71"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("E41").Value = "  +16.14%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.99"
$ws.Range("E42").Value = "  +7.35%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.122"
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.39"
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0355"
$ws.Range("E46").Value = "  +5.96%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.789.17"
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.87"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.77"
$ws.Range("E50").Value = "  +10.26%  "
$ws.Range("E51").Value = "  +4.18%  "
